# Generate Report for Handoff
# The "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d" entry on the zh-cn sheet got a new
# handoff xliff generated, so its "Latest Handoff Datetime" timestamp moves
# forward from 2016-08-24 00:39:40 to 2016-08-24 00:39:56.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("H5").Value = "2016-08-24 00:39:56"
